$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il1b"
$ws.Cells.Item(2, 3).Value = "Il1r2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1413.335253
$ws.Cells.Item(2, 8).Value = 4240.005759
$ws.Cells.Item(2, 9).Value = 0.9999668843963775
$ws.Cells.Item(2, 10).Value = 0.9999668843963775
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 265.842631
$ws.Cells.Item(2, 14).Value = 797.527893
$ws.Cells.Item(2, 15).Value = 0.99055425962745
$ws.Cells.Item(2, 16).Value = 0.99055425962745
$ws.Cells.Item(2, 17).Value = 375724.7621425706
$ws.Cells.Item(2, 18).Value = 3381522.859283135
$ws.Cells.Item(2, 19).Value = 0.9905214568252216
$ws.Cells.Item(2, 20).Value = 0.9905214568252216

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il1b"
$ws.Cells.Item(3, 3).Value = "Il1r2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1413.335253
$ws.Cells.Item(3, 8).Value = 4240.005759
$ws.Cells.Item(3, 9).Value = 0.9999668843963775
$ws.Cells.Item(3, 10).Value = 0.9999668843963775
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.942782333333333
$ws.Cells.Item(3, 14).Value = 5.828347
$ws.Cells.Item(3, 15).Value = 0.007238986871944891
$ws.Cells.Item(3, 16).Value = 0.007238986871944891
$ws.Cells.Item(3, 17).Value = 2745.802760605597
$ws.Cells.Item(3, 18).Value = 24712.22484545037
$ws.Cells.Item(3, 19).Value = 0.007238747148525012
$ws.Cells.Item(3, 20).Value = 0.007238747148525012

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il1b"
$ws.Cells.Item(4, 3).Value = "Il1r2"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1413.335253
$ws.Cells.Item(4, 8).Value = 4240.005759
$ws.Cells.Item(4, 9).Value = 0.9999668843963775
$ws.Cells.Item(4, 10).Value = 0.9999668843963775
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.5922433333333333
$ws.Cells.Item(4, 14).Value = 1.77673
$ws.Cells.Item(4, 15).Value = 0.002206753500604999
$ws.Cells.Item(4, 16).Value = 0.002206753500604999
$ws.Cells.Item(4, 17).Value = 837.03838135423
$ws.Cells.Item(4, 18).Value = 7533.345432188069
$ws.Cells.Item(4, 19).Value = 0.002206680422630781
$ws.Cells.Item(4, 20).Value = 0.002206680422630781

$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Il1b"
$ws.Cells.Item(5, 3).Value = "Il1r2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.04680500000000001
$ws.Cells.Item(5, 8).Value = 0.140415
$ws.Cells.Item(5, 9).Value = 0.00003311560362258399
$ws.Cells.Item(5, 10).Value = 0.00003311560362258399
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 265.842631
$ws.Cells.Item(5, 14).Value = 797.527893
$ws.Cells.Item(5, 15).Value = 0.99055425962745
$ws.Cells.Item(5, 16).Value = 0.99055425962745
$ws.Cells.Item(5, 17).Value = 12.442764343955
$ws.Cells.Item(5, 18).Value = 111.984879095595
$ws.Cells.Item(5, 19).Value = 0.00003280280222848478
$ws.Cells.Item(5, 20).Value = 0.00003280280222848478

$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Il1b"
$ws.Cells.Item(6, 3).Value = "Il1r2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.04680500000000001
$ws.Cells.Item(6, 8).Value = 0.140415
$ws.Cells.Item(6, 9).Value = 0.00003311560362258399
$ws.Cells.Item(6, 10).Value = 0.00003311560362258399
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.942782333333333
$ws.Cells.Item(6, 14).Value = 5.828347
$ws.Cells.Item(6, 15).Value = 0.007238986871944891
$ws.Cells.Item(6, 16).Value = 0.007238986871944891
$ws.Cells.Item(6, 17).Value = 0.09093192711166669
$ws.Cells.Item(6, 18).Value = 0.818387344005
$ws.Cells.Item(6, 19).Value = 0.0000002397234198804162
$ws.Cells.Item(6, 20).Value = 0.0000002397234198804162

$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Il1b"
$ws.Cells.Item(7, 3).Value = "Il1r2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.04680500000000001
$ws.Cells.Item(7, 8).Value = 0.140415
$ws.Cells.Item(7, 9).Value = 0.00003311560362258399
$ws.Cells.Item(7, 10).Value = 0.00003311560362258399
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.5922433333333333
$ws.Cells.Item(7, 14).Value = 1.77673
$ws.Cells.Item(7, 15).Value = 0.002206753500604999
$ws.Cells.Item(7, 16).Value = 0.002206753500604999
$ws.Cells.Item(7, 17).Value = 0.02771994921666667
$ws.Cells.Item(7, 18).Value = 0.24947954295
$ws.Cells.Item(7, 19).Value = 0.00000007307797421878481
$ws.Cells.Item(7, 20).Value = 0.00000007307797421878481

